$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.910.75'
$ws.Range('E2').Value = '  -0.18%  '

$ws.Range('D3').Value = '1.814.44'
$ws.Range('E3').Value = '  +0.23%  '

$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('D5').Value = '309.13'
$ws.Range('E5').Value = '  -0.55%  '

$ws.Range('E6').Value = '  +0.04%  '

$ws.Range('D7').Value = '0.4650'
$ws.Range('E7').Value = '  +0.40%  '

$ws.Range('E8').Value = '  -1.65%  '

$ws.Range('D9').Value = '0.07364'
$ws.Range('E9').Value = '  -0.15%  '

$ws.Range('D10').Value = '0.8685'
$ws.Range('E10').Value = '  -0.77%  '

$ws.Range('D11').Value = '20.23'
$ws.Range('E11').Value = '  -1.17%  '

$ws.Range('D12').Value = '1.832.61'
$ws.Range('E12').Value = '  +0.75%  '

$ws.Range('D13').Value = '5.377'
$ws.Range('E13').Value = '  +0.17%  '

$ws.Range('E14').Value = '  +0.74%  '

$ws.Range('D15').Value = '6.502'
$ws.Range('E15').Value = '  -0.43%  '

$ws.Range('D16').Value = '91.11'
$ws.Range('E16').Value = '  -1.43%  '

$ws.Range('D17').Value = '1.004'
$ws.Range('E17').Value = '  +0.11%  '

$ws.Range('D18').Value = '0.000008687'
$ws.Range('E18').Value = '  -0.34%  '

$ws.Range('E19').Value = '  +0.00%  '

$ws.Range('D20').Value = '14.62'
$ws.Range('E20').Value = '  -0.65%  '

$ws.Range('D21').Value = '26.934.91'
$ws.Range('E21').Value = '  -0.14%  '

$ws.Range('E22').Value = '  -0.51%  '

$ws.Range('D23').Value = '10.57'
$ws.Range('E23').Value = '  -0.78%  '

$ws.Range('D24').Value = '2.064.13'
$ws.Range('E24').Value = '  +0.97%  '

$ws.Range('E25').Value = '  -0.24%  '

$ws.Range('D26').Value = '150.78'
$ws.Range('E26').Value = '  -0.61%  '

$ws.Range('D27').Value = '18.32'
$ws.Range('E27').Value = '  -0.46%  '

$ws.Range('D28').Value = '2.121'
$ws.Range('E28').Value = '  -1.23%  '

$ws.Range('D29').Value = '5.245'
$ws.Range('E29').Value = '  -0.74%  '

$ws.Range('D30').Value = '115.48'
$ws.Range('E30').Value = '  -0.43%  '

$ws.Range('D31').Value = '0.08914'

$ws.Range('D32').Value = '0.7538'
$ws.Range('E32').Value = '  -0.37%  '

$ws.Range('D33').Value = '1.159'
$ws.Range('E33').Value = '  +0.16%  '

$ws.Range('D34').Value = '4.476'
$ws.Range('E34').Value = '  +0.43%  '

$ws.Range('E35').Value = '  -0.88%  '

$ws.Range('E36').Value = '  +0.09%  '

$ws.Range('E37').Value = '  -0.83%  '

$ws.Range('D38').Value = '0.05283'
$ws.Range('E38').Value = '  +0.79%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.985'
$ws.Range('E39').Value = '  +2.20%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.01944'
$ws.Range('E40').Value = '  -1.89%  '

$ws.Range('D41').Value = '7.219'
$ws.Range('E41').Value = '  +0.33%  '

$ws.Range('D42').Value = '0.5285'
$ws.Range('E42').Value = '  -0.55%  '

$ws.Range('D43').Value = '2.298'
$ws.Range('E43').Value = '  -6.27%  '

$ws.Range('D44').Value = '0.1653'
$ws.Range('E44').Value = '  -0.75%  '

$ws.Range('D45').Value = '8.405'
$ws.Range('E45').Value = '  -1.39%  '

$ws.Range('D46').Value = '0.4845'
$ws.Range('E46').Value = '  -2.87%  '

$ws.Range('D47').Value = '10.45'
$ws.Range('E47').Value = '  +0.66%  '

$ws.Range('E48').Value = '  +0.04%  '

$ws.Range('D49').Value = '103.12'
$ws.Range('E49').Value = '  -0.97%  '

$ws.Range('D50').Value = '1.657'
$ws.Range('E50').Value = '  -1.25%  '

$ws.Range('E51').Value = '  -0.13%  '
